$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.387.74"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.570.98"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.61"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3755"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.79"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3408"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.948"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.902"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "1.581.60"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001131"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.34"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06756"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.69"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.206"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5276"
$ws.Range("E23").Value = "  -4.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.97"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "22.383.41"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.730"
$ws.Range("E27").Value = "  -7.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.23"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "144.90"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.067"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.71"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "1.747.14"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.023"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.007"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.145"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.03"
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08550"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02548"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2312"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06500"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.323"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.427"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6428"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.51"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.14"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6010"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.785"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.291"
$ws.Range("E49").Value = "  +8.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.080"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.56"
$ws.Range("E51").Value = "  +3.84%  "
